$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be pinned to Text format first,
# otherwise Excel auto-converts them to numbers (losing the trailing/format digits),
# then the style is reset back to Normal so no stray formatting is left behind.
$textCells = @(
  'D5',
  'D6',
  'D7',
  'D9',
  'D10',
  'D11',
  'D12',
  'D14',
  'D15',
  'D16',
  'D17',
  'D21',
  'D22',
  'D23',
  'D27',
  'D28',
  'D29',
  'D30',
  'D33',
  'D35',
  'D36',
  'D37',
  'D38',
  'D39',
  'D41',
  'D43',
  'D44',
  'D45',
  'D46',
  'D47',
  'D48'
)
foreach ($c in $textCells) {
  $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.316.90'
$ws.Range('E2').Value = '  -1.49%  '

$ws.Range('D3').Value = '2.339.06'
$ws.Range('E3').Value = '  +3.06%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '0.650'
$ws.Range('E5').Value = '  +3.29%  '

$ws.Range('D6').Value = '231.89'
$ws.Range('E6').Value = '  +0.31%  '

$ws.Range('D7').Value = '65.71'
$ws.Range('E7').Value = '  +3.17%  '

$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('D9').Value = '0.452'
$ws.Range('E9').Value = '  +1.08%  '

$ws.Range('D10').Value = '0.0953'
$ws.Range('E10').Value = '  -4.18%  '

$ws.Range('D11').Value = '56.90'
$ws.Range('E11').Value = '  -1.01%  '

$ws.Range('D12').Value = '26.72'
$ws.Range('E12').Value = '  -2.19%  '

$ws.Range('D13').Value = '2.688.10'
$ws.Range('E13').Value = '  +3.03%  '

$ws.Range('D14').Value = '0.105'
$ws.Range('E14').Value = '  -1.50%  '

$ws.Range('D15').Value = '15.33'
$ws.Range('E15').Value = '  -2.93%  '

$ws.Range('D16').Value = '6.22'
$ws.Range('E16').Value = '  +1.47%  '

$ws.Range('D17').Value = '0.843'
$ws.Range('E17').Value = '  +0.70%  '

$ws.Range('D18').Value = '2.337.02'
$ws.Range('E18').Value = '  +2.54%  '

$ws.Range('D19').Value = '43.249.89'
$ws.Range('E19').Value = '  -1.54%  '

$ws.Range('D20').Value = '0.0₃0975'
$ws.Range('E20').Value = '  -3.21%  '

$ws.Range('D21').Value = '73.99'
$ws.Range('E21').Value = '  +0.23%  '

$ws.Range('D22').Value = '6.18'
$ws.Range('E22').Value = '  +1.17%  '

$ws.Range('D23').Value = '248.54'
$ws.Range('E23').Value = '  -1.62%  '

$ws.Range('E24').Value = '  +16.38%  '

$ws.Range('E25').Value = '  -0.04%  '

$ws.Range('E26').Value = '  -1.18%  '

$ws.Range('D27').Value = '2.22'
$ws.Range('E27').Value = '  -1.54%  '

$ws.Range('D28').Value = '9.88'
$ws.Range('E28').Value = '  -2.29%  '

$ws.Range('D29').Value = '175.47'
$ws.Range('E29').Value = '  +2.34%  '

$ws.Range('D30').Value = '22.14'
$ws.Range('E30').Value = '  +5.80%  '

$ws.Range('E31').Value = '  +6.13%  '

$ws.Range('E32').Value = '  -7.94%  '

$ws.Range('D33').Value = '0.125'
$ws.Range('E33').Value = '  +0.25%  '

$ws.Range('E34').Value = '  +4.02%  '

$ws.Range('D35').Value = '0.0686'
$ws.Range('E35').Value = '  -2.86%  '

$ws.Range('D36').Value = '4.93'
$ws.Range('E36').Value = '  +1.46%  '

$ws.Range('D37').Value = '2.51'
$ws.Range('E37').Value = '  +7.98%  '

$ws.Range('D38').Value = '6.45'
$ws.Range('E38').Value = '  -0.82%  '

$ws.Range('D39').Value = '3.59'
$ws.Range('E39').Value = '  -5.94%  '

$ws.Range('E40').Value = '  -3.34%  '

$ws.Range('D41').Value = '8.99'
$ws.Range('E41').Value = '  +8.62%  '

$ws.Range('E42').Value = '  +0.18%  '

$ws.Range('D43').Value = '18.00'
$ws.Range('E43').Value = '  +2.48%  '

$ws.Range('D44').Value = '1.17'
$ws.Range('E44').Value = '  +7.98%  '

$ws.Range('D45').Value = '98.86'
$ws.Range('E45').Value = '  +0.62%  '

$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').Value = '4.42'
$ws.Range('E46').Value = '  +0.67%  '

$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Value = '1.19'
$ws.Range('E47').Value = '  -1.09%  '

$ws.Range('D48').Value = '0.0944'
$ws.Range('E48').Value = '  -4.53%  '

$ws.Range('D49').Value = '1.435.79'
$ws.Range('E49').Value = '  -0.75%  '

$ws.Range('E50').Value = '  -6.43%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.562.80'
$ws.Range('E51').Value = '  +3.14%  '

foreach ($c in $textCells) {
  $ws.Range($c).Style = "Normal"
}
